$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E6").Value = "['Normal']"

$ws.Range("D16").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E16").Value = "['Normal']"

$ws.Range("D30").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E30").Value = "['Normal']"

$ws.Range("D36").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E36").Value = "['Normal', 'HardwareFault']"
